# "proper simulation for one team"
# The 2014-15 ("14-15") sheet mixed full team names (row labels, column A)
# with abbreviation codes (column headers, row 1). Normalize column A to use
# the same abbreviation codes as the header row, and fix a handful of
# abbreviations/names that were stale for that season (BRK->NJN, CHO->CHA,
# NOP->NOH, OKC->SEA), matching the "09-10" sheet's already-correct coding.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("14-15")

# Column order (same list used for both the header row and the row labels).
$teams = @("ATL","BOS","NJN","CHI","CHA","CLE","DAL","DEN","DET","GSW","HOU","IND","LAC","LAL","MEM","MIA","MIL","MIN","NOH","NYK","SEA","ORL","PHI","PHO","POR","SAC","SAS","TOR","UTA","WAS")

# Header row (B1:AE1) - fix the four stale codes.
for ($i = 0; $i -lt $teams.Length; $i++) {
    $ws1.Cells.Item(1, $i + 2).Value = $teams[$i]
}

# Row labels (A2:A31) - switch from full franchise names to the same
# abbreviation codes used in the header.
for ($i = 0; $i -lt $teams.Length; $i++) {
    $ws1.Cells.Item($i + 2, 1).Value = $teams[$i]
}

# Column A is now short codes instead of long franchise names - narrow it
# back down to a sensible display width.
$ws1.Columns.Item(1).ColumnWidth = 14

# Make the "14-15" sheet the active tab/selection, leaving "09-10" as it was.
$ws1.Activate()
$ws1.Range("A11").Select()
